$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 365, pushing existing rows 365-387 down to 366-388.
$ws.Rows.Item(365).Insert()

# Populate the new row 365 with the new record.
$ws.Cells.Item(365, 1).Value = 7
$ws.Cells.Item(365, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(365, 3).Value = "Ñuble"
$ws.Cells.Item(365, 4).Value = 45265
$ws.Cells.Item(365, 5).Value = 16
$ws.Cells.Item(365, 6).Value = 100112032
$ws.Cells.Item(365, 7).Value = "Zapallo italiano"
$ws.Cells.Item(365, 8).Value = "Sin especificar"
$ws.Cells.Item(365, 9).Value = "Primera"
$ws.Cells.Item(365, 10).Value = 60
$ws.Cells.Item(365, 11).Value = 12000
$ws.Cells.Item(365, 12).Value = 12000
$ws.Cells.Item(365, 13).Value = 12000
$ws.Cells.Item(365, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(365, 15).Value = "Región del Maule"
$ws.Cells.Item(365, 16).Value = 240
$ws.Cells.Item(365, 17).Value = 50
$ws.Cells.Item(365, 18).Value = "Hortaliza"
